# Atualização de bases das ligas, do dia: 30-05-2024 às 23:16
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns I (ht_goals_h -> HTHG) and J (ht_goals_a -> HTAG)
$ws.Cells.Item(1, 9).Value = "HTHG"
$ws.Cells.Item(1, 10).Value = "HTAG"

# Re-write rows 195-197 (match id 193-195) with the updated match data.
# Row 195 (id 193)
$ws.Cells.Item(195, 2).Value = 8209692
$ws.Cells.Item(195, 5).Value = "Puskas Academy"
$ws.Cells.Item(195, 6).Value = "Debreceni VSC"
$ws.Cells.Item(195, 7).Value = 4
$ws.Cells.Item(195, 8).Value = 1
$ws.Cells.Item(195, 9).Value = 2
$ws.Cells.Item(195, 10).Value = 1
$ws.Cells.Item(195, 11).Value = "H"
$ws.Cells.Item(195, 12).Value = 1.7
$ws.Cells.Item(195, 13).Value = 3.6
$ws.Cells.Item(195, 14).Value = 4.5
$ws.Cells.Item(195, 15).Value = 1.333
$ws.Cells.Item(195, 16).Value = 4.5
$ws.Cells.Item(195, 17).Value = 8
$ws.Cells.Item(195, 18).Value = -1.5
$ws.Cells.Item(195, 19).Value = 2.025
$ws.Cells.Item(195, 20).Value = 1.825
$ws.Cells.Item(195, 21).Value = 3
$ws.Cells.Item(195, 22).Value = 1.875
$ws.Cells.Item(195, 23).Value = 1.975
$ws.Cells.Item(195, 24).Value = 0.333
$ws.Cells.Item(195, 25).Value = -1
$ws.Cells.Item(195, 26).Value = -1
$ws.Cells.Item(195, 27).Value = 1.025
$ws.Cells.Item(195, 28).Value = -1
$ws.Cells.Item(195, 29).Value = 0.875
$ws.Cells.Item(195, 30).Value = -1

# Row 196 (id 194)
$ws.Cells.Item(196, 2).Value = 8209693
$ws.Cells.Item(196, 5).Value = "MOL Fehervar FC"
$ws.Cells.Item(196, 6).Value = "Diosgyori VTK"
$ws.Cells.Item(196, 7).Value = 0
$ws.Cells.Item(196, 8).Value = 0
$ws.Cells.Item(196, 9).Value = 0
$ws.Cells.Item(196, 10).Value = 0
$ws.Cells.Item(196, 11).Value = "D"
$ws.Cells.Item(196, 12).Value = 1.571
$ws.Cells.Item(196, 13).Value = 4
$ws.Cells.Item(196, 14).Value = 5
$ws.Cells.Item(196, 15).Value = 1.6
$ws.Cells.Item(196, 16).Value = 3.9
$ws.Cells.Item(196, 17).Value = 4.75
$ws.Cells.Item(196, 18).Value = -1
$ws.Cells.Item(196, 19).Value = 2.025
$ws.Cells.Item(196, 20).Value = 1.825
$ws.Cells.Item(196, 21).Value = 3
$ws.Cells.Item(196, 22).Value = 1.925
$ws.Cells.Item(196, 23).Value = 1.925
$ws.Cells.Item(196, 24).Value = -1
$ws.Cells.Item(196, 25).Value = 2.9
$ws.Cells.Item(196, 26).Value = -1
$ws.Cells.Item(196, 27).Value = -1
$ws.Cells.Item(196, 28).Value = 0.825
$ws.Cells.Item(196, 29).Value = -1
$ws.Cells.Item(196, 30).Value = 0.925

# Row 197 (id 195)
$ws.Cells.Item(197, 2).Value = 8209690
$ws.Cells.Item(197, 5).Value = "Paksi"
$ws.Cells.Item(197, 6).Value = "Kisvarda FC"
$ws.Cells.Item(197, 7).Value = 2
$ws.Cells.Item(197, 8).Value = 1
$ws.Cells.Item(197, 9).Value = 1
$ws.Cells.Item(197, 10).Value = 0
$ws.Cells.Item(197, 11).Value = "H"
$ws.Cells.Item(197, 12).Value = 1.444
$ws.Cells.Item(197, 13).Value = 4.333
$ws.Cells.Item(197, 14).Value = 6
$ws.Cells.Item(197, 15).Value = 1.45
$ws.Cells.Item(197, 16).Value = 4.75
$ws.Cells.Item(197, 17).Value = 5
$ws.Cells.Item(197, 18).Value = -1.25
$ws.Cells.Item(197, 19).Value = 2.025
$ws.Cells.Item(197, 20).Value = 1.825
$ws.Cells.Item(197, 21).Value = 3.25
$ws.Cells.Item(197, 22).Value = 2
$ws.Cells.Item(197, 23).Value = 1.85
$ws.Cells.Item(197, 24).Value = 0.45
$ws.Cells.Item(197, 25).Value = -1
$ws.Cells.Item(197, 26).Value = -1
$ws.Cells.Item(197, 27).Value = -0.5
$ws.Cells.Item(197, 28).Value = 0.4125
$ws.Cells.Item(197, 29).Value = -0.5
$ws.Cells.Item(197, 30).Value = 0.425
